$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 28
$ws.Range("E11").Value = 14
$ws.Range("E15").Value = 91
$ws.Range("F15").Value = 44
$ws.Range("H15").Value = 44
$ws.Range("E16").Value = 303
